# "rename best rolling to critical"
#
# gc_fields_uom: append 4 new unit-definition rows for the
# "__CalcBestRollingWeightedMean*" fields at the end of the table.
#
# gc_fields_display: insert 4 new display-name rows for the same fields,
# in their correct alphabetical slot (right after __CalcAscentSpeed),
# pushing the existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: gc_fields_uom  (A1:D129 -> A1:D133)
# ---------------------------------------------------------------
$wsUom = $wb.Worksheets.Item("gc_fields_uom")

$uomRows = @(
    @("__CalcBestRollingWeightedMeanHeartRate", "all", "bpm",      "bpm"),
    @("__CalcBestRollingWeightedMeanPace",      "all", "minperkm", "minpermile"),
    @("__CalcBestRollingWeightedMeanPower",     "all", "watt",     "watt"),
    @("__CalcBestRollingWeightedMeanSpeed",     "all", "kph",      "mph")
)

$startRow = 130
for ($i = 0; $i -lt $uomRows.Length; $i++) {
    $r = $startRow + $i
    $row = $uomRows[$i]
    $wsUom.Cells.Item($r, 1).Value = $row[0]
    $wsUom.Cells.Item($r, 2).Value = $row[1]
    $wsUom.Cells.Item($r, 3).Value = $row[2]
    $wsUom.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------
# Sheet: gc_fields_display  (A1:I195 -> A1:I199)
# ---------------------------------------------------------------
$wsDisplay = $wb.Worksheets.Item("gc_fields_display")

# Insert 4 blank rows right before the existing __CalcDescentSpeed row
# (row 153), pushing it (and everything after) down by 4.
$wsDisplay.Range("A153:A156").EntireRow.Insert()

$displayRows = @(
    @("__CalcBestRollingWeightedMeanHeartRate", "Critical Heart Rate"),
    @("__CalcBestRollingWeightedMeanPace",      "Critical Pace"),
    @("__CalcBestRollingWeightedMeanPower",     "Critical Power"),
    @("__CalcBestRollingWeightedMeanSpeed",     "Critical Speed")
)

$startRow = 153
for ($i = 0; $i -lt $displayRows.Length; $i++) {
    $r = $startRow + $i
    $row = $displayRows[$i]
    $wsDisplay.Cells.Item($r, 1).Value = $row[0]
    $wsDisplay.Cells.Item($r, 2).Value = $row[1]
}
